$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Excel alignment constants ---
$xlVAlignTop = -4160

# Helper: reset a cell to the plain "Normal" style first (avoids the engine
# stacking a spurious explicit-but-empty alignment record on top of
# whatever formatting the cell previously had), then re-apply only the
# formatting actually needed.
function Set-CellFormat($cell, [bool]$wrap, [bool]$top) {
    $cell.Style = "Normal"
    if ($wrap) { $cell.WrapText = $true }
    if ($top) { $cell.VerticalAlignment = $xlVAlignTop }
}

# --- 1. Snapshot current (pre-edit) row 3 and row 4 contents/formatting (B..E) ---
$row3 = @{}
$row4 = @{}
foreach ($c in 2..5) {
    $cell3 = $ws.Cells.Item(3, $c)
    $row3[$c] = @{ Value = $cell3.Value2; Wrap = [bool]$cell3.WrapText; Top = ($cell3.VerticalAlignment -eq $xlVAlignTop) }

    $cell4 = $ws.Cells.Item(4, $c)
    $row4[$c] = @{ Value = $cell4.Value2; Wrap = [bool]$cell4.WrapText; Top = ($cell4.VerticalAlignment -eq $xlVAlignTop) }
}

# --- 2. Write old row4 content/formatting into row3 (B..E) ---
foreach ($c in 2..5) {
    $cell = $ws.Cells.Item(3, $c)
    $src = $row4[$c]
    $cell.Value2 = $src.Value
    if ($c -ne 4) {
        # D column keeps its Hyperlink style (s=6) untouched; only B,C,E
        # need their wrap/vertical-alignment re-derived from the row that
        # is moving in.
        Set-CellFormat $cell $src.Wrap $src.Top
    }
}

# --- 3. Write old row3 content/formatting into row4 (B..E) ---
foreach ($c in 2..5) {
    $cell = $ws.Cells.Item(4, $c)
    $src = $row3[$c]
    $cell.Value2 = $src.Value
    if ($c -ne 4) {
        Set-CellFormat $cell $src.Wrap $src.Top
    }
}

# --- 4. Update the "verify ..." task text that now lives in E4 ---
$ws.Cells.Item(4, 5).Value2 = 'verify text equals "HART, WILLIAM S." in "patient_name"'

# --- 5. Re-number the id column to stay sequential ---
$ws.Cells.Item(4, 1).Value2 = 3
$ws.Cells.Item(5, 1).Value2 = 4

# --- 6. Row heights follow the content that now occupies each row ---
$ws.Rows.Item(3).RowHeight = 78.75
$ws.Rows.Item(4).RowHeight = 31.5

# --- 7. Hyperlink "location" (sub-address) moves from D4 to D3 ---
$hlD4 = $null
$hlD3 = $null
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$D`$4") { $hlD4 = $hl }
    if ($hl.Range.Address() -eq "`$D`$3") { $hlD3 = $hl }
}
if ($hlD4 -ne $null) { $hlD4.SubAddress = "" }
if ($hlD3 -ne $null) { $hlD3.SubAddress = "/home" }

# --- 8. Selection moves from E3 to A4 ---
$ws.Range("A4").Select() | Out-Null
